$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 currently holds a "workingFutures" placeholder with only a few
# cells populated (metric, subheading, dataText, mapComment). We are
# replacing it with the full "wfEmployment" (Working Futures forecasted
# employment) data row, matching the format used by the other metric rows.

# Copy the formatting from row 2 (same style pattern as the target row)
# onto row 21 first, so every cell (A21:M21) carries the right style index.
$ws.Range("A2:M2").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now populate the new values for row 21.
$ws.Range("A21").Value = "wfEmployment"
$ws.Range("B21").Value = "2035"
$ws.Range("C21").Value = "To come"
$ws.Range("D21").Value = "To come"
$ws.Range("E21").Value = "To come"
$ws.Range("F21").Value = "To come"
$ws.Range("G21").Value = "Forecasted employment"
$ws.Range("H21").Value = "forecasted employment volume growth"
$ws.Range("I21").Value = "Forecasted employment growth in"
$ws.Range("J21").Value = "Forecasted employment"
$ws.Range("K21").Value = "forecasted employment volumes"
$ws.Range("L21").Value = "forecasted employment volume"
$ws.Range("M21").Value = "Forecasted employment is"

# Match the row height used for the new content.
$ws.Rows("21").RowHeight = 42.5

# Update the saved view state (scroll position / active selection) to
# reflect where the author left the sheet after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 20
$win.ScrollColumn = 1
$ws.Range("I22").Select()
